$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

$newText = @'
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 4.57 = 18059.36 pesos
✅ 18059.36 pesos = 4.53 = 937.62 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
'@

$ws1.Range("A1").Value = $newText

$ws2.Range("N10").Value = 219
$ws2.Range("O10").Value = 3955
$ws2.Range("N12").Value = 3987.96
